# DEV 6 - CLI Changes for Manager Application and Registration
#
# Row 2 (existing Officer Registration #1): status goes from "Approved" -> "Pending"
# Row 3 (existing entry): becomes Officer Registration #2, Project ID 1,
#                          status "Successful" (NRIC/date unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Registration Status Approved -> Pending
$ws.Range("D2").Value = "Pending"

# Row 3: Officer Registration ID, Project ID, Registration Status updates
$ws.Range("A3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Successful"

# Move the active cell selection to D4 (matches the saved view state)
$ws.Range("D4").Select()
